# Re-generate the statistics with fixed minutes/seconds zero-padding in the
# "haul" (Общее время) time fields, e.g. "317 ч. 55 мин. 2 сек." ->
# "317 ч. 55 мин. 02 сек.".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$numRows = $used.Rows.Count
$numCols = $used.Columns.Count

$re = [regex]'^(\d+) ч\. (\d+) мин\. (\d+) сек\.$'

$changed = 0

for ($r = 1; $r -le $numRows; $r++) {
    for ($c = 1; $c -le $numCols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value2
        if ($v -eq $null) { continue }
        if ($v.GetType().Name -ne "String") { continue }

        $m = $re.Match($v)
        if (-not $m.Success) { continue }

        $hours = $m.Groups[1].Value
        $minutes = $m.Groups[2].Value
        $seconds = $m.Groups[3].Value

        # Avoid string "-eq"/"-ne" on numeric-looking strings: this engine's
        # PowerShell coerces both sides to numbers for comparison (so "02"
        # -eq "2" is True). Compare .Length (a real int) instead.
        $needsPad = ($minutes.Length -lt 2) -or ($seconds.Length -lt 2)
        if (-not $needsPad) { continue }

        $paddedMinutes = $minutes.PadLeft(2, '0')
        $paddedSeconds = $seconds.PadLeft(2, '0')

        $newValue = "$hours ч. $paddedMinutes мин. $paddedSeconds сек."
        $cell.Value2 = $newValue
        $changed = $changed + 1
    }
}

Write-Host "Cells updated: $changed"
